$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row => New C (nombre_aides) value, New D (montant_total) value
$updates = @(
    @(2, 33904, 49175127),
    @(3, 83038, 121991959),
    @(4, 28414, 42143697),
    @(5, 7743, 11521714),
    @(6, 1587, 2365128),
    @(7, 112, 163093),
    @(10, 36925, 50396966),
    @(11, 8743, 12672550),
    @(12, 23957, 35204496),
    @(13, 7606, 11308298),
    @(14, 1906, 2835211),
    @(15, 337, 495964),
    @(18, 9205, 12264409),
    @(19, 12253, 17724167),
    @(20, 29306, 43097864),
    @(21, 9494, 14128387),
    @(22, 2369, 3529219),
    @(23, 386, 576342),
    @(25, 10599, 14279971),
    @(26, 6966, 10109693),
    @(27, 20744, 30500419),
    @(28, 7206, 10729039),
    @(29, 1762, 2636482),
    @(30, 270, 402415),
    @(32, 7497, 9970944),
    @(33, 2699, 3897059),
    @(34, 6809, 9957833),
    @(35, 2760, 4086732),
    @(36, 730, 1087763),
    @(37, 130, 194304),
    @(39, 2020, 2726045),
    @(40, 15551, 22539995),
    @(41, 46822, 68741146),
    @(42, 17499, 26005445),
    @(43, 5011, 7469155),
    @(44, 923, 1379950),
    @(48, 14954, 20042130),
    @(49, 1599, 2319448),
    @(50, 5654, 8326804),
    @(51, 1993, 2977957),
    @(52, 644, 962026),
    @(55, 4731, 6555057),
    @(56, 609, 891985),
    @(57, 1559, 2309923),
    @(58, 607, 906419),
    @(62, 884, 1250750),
    @(63, 13915, 20129217),
    @(64, 41343, 60605710),
    @(65, 14537, 21631516),
    @(66, 4124, 6148083),
    @(67, 758, 1130520),
    @(68, 64, 94689),
    @(70, 13854, 18365019),
    @(71, 3946, 5732057),
    @(72, 10286, 15124381),
    @(73, 3519, 5246526),
    @(74, 1199, 1793091),
    @(78, 4570, 6170443),
    @(79, 1313, 1897497),
    @(80, 4381, 6455062),
    @(81, 1706, 2546932),
    @(82, 599, 897641),
    @(83, 141, 211069),
    @(86, 2839, 3778957),
    @(87, 488, 728159),
    @(90, 32, 48000),
    @(92, 9774, 14229434),
    @(93, 27088, 39867767),
    @(94, 9059, 13481612),
    @(95, 2444, 3647399),
    @(96, 405, 605249),
    @(99, 8875, 11786200),
    @(100, 27511, 39754237),
    @(101, 61094, 89559023),
    @(102, 19749, 29384531),
    @(103, 5445, 8118982),
    @(104, 938, 1401893),
    @(105, 51, 73358),
    @(108, 23389, 31398223),
    @(109, 31996, 46269732),
    @(110, 69701, 102117634),
    @(111, 21791, 32369780),
    @(112, 5704, 8486863),
    @(113, 977, 1457677),
    @(114, 47, 68728),
    @(115, 15, 22500),
    @(117, 28197, 37636643),
    @(118, 12045, 17474245),
    @(119, 30028, 44165716),
    @(120, 10650, 15837716),
    @(121, 2674, 3989370),
    @(122, 394, 584990),
    @(125, 9874, 13268415),
    @(126, 31280, 45235774),
    @(127, 73956, 108495538),
    @(128, 22267, 33125553),
    @(129, 5712, 8529097),
    @(130, 1173, 1749220),
    @(133, 26275, 35684768),
    @(134, 43838, 63898294),
    @(135, 127780, 188535363),
    @(136, 56098, 83654110),
    @(137, 17705, 26462934),
    @(138, 3756, 5613773),
    @(139, 201, 299355),
    @(140, 18, 25905),
    @(145, 43448, 59651182),
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 3).Value = $u[1]
    $ws.Cells.Item($row, 4).Value = $u[2]
}
